$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 is a brand new row that receives the ORIGINAL values that used
# to live in row 39 (before this edit).
$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value = 44442
$ws.Cells.Item(40, 4).NumberFormat = $ws.Cells.Item(39, 4).NumberFormat()
$ws.Cells.Item(40, 5).Value = 15
$ws.Cells.Item(40, 6).Value = 100112052
$ws.Cells.Item(40, 7).Value = "Albahaca"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 240
$ws.Cells.Item(40, 11).Value = 2300
$ws.Cells.Item(40, 12).Value = 2500
$ws.Cells.Item(40, 13).Value = 2400
$ws.Cells.Item(40, 14).Value = "$/paquete"
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 2400
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"

# Row 39 is updated in place with the new weekly figures.
$ws.Cells.Item(39, 4).Value = 44747
$ws.Cells.Item(39, 10).Value = 250
$ws.Cells.Item(39, 11).Value = 2000
$ws.Cells.Item(39, 12).Value = 2500
$ws.Cells.Item(39, 13).Value = 2250
$ws.Cells.Item(39, 16).Value = 2250
